$d = $word.ActiveDocument

# 1. "Oude KMSK Navbar: " - merge runs split around "Navbar" (proofErr spell-check marks) into one run.
$d.Content.Find.Execute("Oude KMSK Navbar: ", $false, $false, $false, $false, $false, $true, 1, $false, "Oude KMSK Navbar: ", 2)

# 2. Long nav-bar text line - merge runs split around "Snelschaak" / "k.schap" / "Elo" into one run.
$d.Content.Find.Execute("Home – Bestuur - Ligging lokaal – Kalender – Stapjestornooi – Jeugdkalender – Jeugdschaak – Clubkampioenschap – Laddertornooi – Interclub - Zilveren Toren - Snelschaak k.schap - Trofee Walter Huyck – Wintertornooi – Verzekering – Elo – Historiek – Links - Inlog bestuur", $false, $false, $false, $false, $false, $true, 1, $false, "Home – Bestuur - Ligging lokaal – Kalender – Stapjestornooi – Jeugdkalender – Jeugdschaak – Clubkampioenschap – Laddertornooi – Interclub - Zilveren Toren - Snelschaak k.schap - Trofee Walter Huyck – Wintertornooi – Verzekering – Elo – Historiek – Links - Inlog bestuur", 2)

# 3. "Nieuw KMSK Navbar:" - merge runs split around "Navbar" into one run.
$d.Content.Find.Execute("Nieuw KMSK Navbar:", $false, $false, $false, $false, $false, $true, 1, $false, "Nieuw KMSK Navbar:", 2)

# 4. Content change: drop the "Links (Elo – Links - Historiek)" expansion, add "- Kalender" to Team Competitie,
#    and collapse "Links" down to a plain link (no parenthetical).
$d.Content.Find.Execute("Team Competitie (Interclub – Zilveren Toren) – Links (Elo – Links - Historiek) – Jeugdschaak", $false, $false, $false, $false, $false, $true, 1, $false, "Team Competitie (Interclub – Zilveren Toren - Kalender) – Links – Jeugdschaak", 2)

# 5. "HomePage: " - merge runs split around "HomePage" (proofErr spell-check marks) into one run.
#    The paragraph starts with the proofErr marker (nothing before it), so a plain Find/Replace
#    can't absorb it (it only clears proofErr nodes sitting *between* two runs it rewrites).
#    Insert a zero-width marker run just before it first, then replace across that boundary so the
#    proofErr ends up "between" runs being merged and gets dropped; the temporary marker is removed
#    as part of the same replace.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "HomePage*") {
        $marker = $d.Range($para.Range.Start, $para.Range.Start)
        $marker.InsertBefore("@")
        break
    }
}
$d.Content.Find.Execute("@HomePage: ", $false, $false, $false, $false, $false, $true, 1, $false, "HomePage: ", 2)
